# Function Efficiency formula update
# Functional if voltage input and output of DC/DC are identical to the
# associated formula.

$wb = $excel.ActiveWorkbook

$wsDCDC     = $wb.Worksheets.Item("DCDC")
$wsCONSUMER = $wb.Worksheets.Item("CONSUMER")

# --- DCDC: fix the efficiency formula text -----------------------------
# The decimal separator in the split threshold was a comma, which breaks
# parsing when the function is evaluated; correct it to a period.
$wsDCDC.Range("B9").Value = "12/5/0.1/-669890*x**4+176938*x**3-16759*x**2+689*x+77/1*x**3-9*x**2+14*x+87"

# widen column B on DCDC so the updated formula text remains legible
$wsDCDC.Columns.Item(2).ColumnWidth = 27.42578125

# --- CONSUMER: add a new "Test" device in column D ----------------------
$wsCONSUMER.Range("D1").Value = "Test"
$wsCONSUMER.Range("D2").Value = "sdf"
$wsCONSUMER.Range("D3").Value = "sdf"
$wsCONSUMER.Range("D4").Value = "sdf"
$wsCONSUMER.Range("D5").Value = 5
$wsCONSUMER.Range("D6").Value = 200

# --- view/selection state ------------------------------------------------
$wsDCDC.Activate() | Out-Null
$wsDCDC.Application.ActiveWindow.Zoom = 175
$wsDCDC.Range("D10").Select() | Out-Null

$wsCONSUMER.Activate() | Out-Null
$wsCONSUMER.Range("D6").Select() | Out-Null
